# "Bao cao buoi 1" slide 1 fix:
# The "TextBox 4" shape holds the title line "Dự án   : " where the
# colon was preceded by three spaces left over from an earlier edit.
# Tidy it up to "Dự án: " (single space before the colon, none before
# it besides the word break already supplied by the previous run).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(4)          # "TextBox 4"
$tr = $sh.TextFrame.TextRange

# Replace only the stray "   : " run with ": " - leaves every other
# run (and its formatting) in the paragraph untouched.
$tr.Replace("   : ", ": ") | Out-Null
